# Applies the change described in the commit:
# "removed the calculations that are not included in the paper;
#  updated tables to include confidence intervals"
#
# On the "summary" worksheet, add two new columns (D, E) with the
# lower/upper bound of the 95% confidence interval for the average AUC,
# computed from the existing Average/StdDev-of-auc columns (B, C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")

# Headers
$ws.Range("D1").Value = "lower interval 95%"
$ws.Range("E1").Value = "upper interval 95%"

# Row 2 (ag_sononet)
$ws.Range("D2").Formula = "=B2-C2/SQRT(5)*1.96"
$ws.Range("E2").Formula = "=B2+C2/SQRT(5)*1.96"

# Row 3 (sononet)
$ws.Range("D3").Formula = "=B3-C3/SQRT(5)*1.96"
$ws.Range("E3").Formula = "=B3+C3/SQRT(5)*1.96"

# Row 4 (Grand Total) intentionally left without a confidence interval,
# matching the source workbook.
